# Auto-generated edit script for cryptos.xlsx update (crypto price refresh).
# - Updates Price (col D) and Volume(1h) (col E) for most coin rows.
# - Row 35/36: swaps Hedera and InjectiveProtocol (with their link/price/volume).
# - Row 49/50: swaps BinanceUSD and ARBITRUM (with their link/price/volume).
# - Row 51: TrustWalletToken replaced by Aave (with its link/price/volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.736.67'
$ws.Cells.Item(2, 5).Value = '  +3.18%  '
$ws.Cells.Item(3, 4).Value = '2.292.67'
$ws.Cells.Item(3, 5).Value = '  +4.85%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(5, 4).Value = '''251.17'
$ws.Cells.Item(5, 5).Value = '  +0.55%  '
$ws.Cells.Item(6, 4).Value = '''0.634'
$ws.Cells.Item(6, 5).Value = '  +2.55%  '
$ws.Cells.Item(7, 4).Value = '''72.73'
$ws.Cells.Item(7, 5).Value = '  +8.16%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 4).Value = '''0.647'
$ws.Cells.Item(9, 5).Value = '  +5.00%  '
$ws.Cells.Item(10, 4).Value = '''39.25'
$ws.Cells.Item(10, 5).Value = '  +1.57%  '
$ws.Cells.Item(11, 5).Value = '  +3.47%  '
$ws.Cells.Item(12, 4).Value = '''59.18'
$ws.Cells.Item(12, 5).Value = '  -0.55%  '
$ws.Cells.Item(13, 4).Value = '''7.31'
$ws.Cells.Item(13, 5).Value = '  +4.43%  '
$ws.Cells.Item(14, 5).Value = '  +1.62%  '
$ws.Cells.Item(15, 4).Value = '2.631.72'
$ws.Cells.Item(15, 5).Value = '  +4.67%  '
$ws.Cells.Item(16, 4).Value = '''15.18'
$ws.Cells.Item(16, 5).Value = '  +4.66%  '
$ws.Cells.Item(17, 4).Value = '''0.884'
$ws.Cells.Item(17, 5).Value = '  +3.34%  '
$ws.Cells.Item(18, 4).Value = '2.272.62'
$ws.Cells.Item(18, 5).Value = '  +3.86%  '
$ws.Cells.Item(19, 4).Value = '42.688.78'
$ws.Cells.Item(19, 5).Value = '  +3.34%  '
$ws.Cells.Item(20, 4).Value = '''0.0' + [char]0x2082 + '01000'
$ws.Cells.Item(20, 5).Value = '  +5.20%  '
$ws.Cells.Item(21, 5).Value = '  +3.46%  '
$ws.Cells.Item(22, 4).Value = '''72.71'
$ws.Cells.Item(22, 5).Value = '  +1.11%  '
$ws.Cells.Item(23, 5).Value = '  +12.94%  '
$ws.Cells.Item(24, 4).Value = '''235.58'
$ws.Cells.Item(24, 5).Value = '  +2.40%  '
$ws.Cells.Item(25, 5).Value = '  +1.79%  '
$ws.Cells.Item(26, 4).Value = '''11.66'
$ws.Cells.Item(26, 5).Value = '  +3.25%  '
$ws.Cells.Item(27, 5).Value = '  -0.04%  '
$ws.Cells.Item(28, 4).Value = '''2.44'
$ws.Cells.Item(28, 5).Value = '  +1.60%  '
$ws.Cells.Item(29, 4).Value = '''3.66'
$ws.Cells.Item(29, 5).Value = '  -0.34%  '
$ws.Cells.Item(30, 5).Value = '  +5.75%  '
$ws.Cells.Item(31, 4).Value = '''167.59'
$ws.Cells.Item(31, 5).Value = '  +0.38%  '
$ws.Cells.Item(32, 4).Value = '''21.14'
$ws.Cells.Item(32, 5).Value = '  +4.45%  '
$ws.Cells.Item(33, 4).Value = '''6.44'
$ws.Cells.Item(33, 5).Value = '  +10.51%  '
$ws.Cells.Item(34, 5).Value = '  +7.31%  '
$ws.Cells.Item(35, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(35, 4).Value = '''31.97'
$ws.Cells.Item(35, 5).Value = '  +24.20%  '
$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).Value = '''0.0806'
$ws.Cells.Item(36, 5).Value = '  +2.36%  '
$ws.Cells.Item(37, 5).Value = '  +3.42%  '
$ws.Cells.Item(38, 4).Value = '''4.78'
$ws.Cells.Item(38, 5).Value = '  +14.85%  '
$ws.Cells.Item(39, 4).Value = '''4.76'
$ws.Cells.Item(39, 5).Value = '  +4.47%  '
$ws.Cells.Item(40, 4).Value = '''0.0308'
$ws.Cells.Item(40, 5).Value = '  +0.75%  '
$ws.Cells.Item(41, 5).Value = '  +19.88%  '
$ws.Cells.Item(42, 5).Value = '  +6.51%  '
$ws.Cells.Item(43, 4).Value = '''6.02'
$ws.Cells.Item(43, 5).Value = '  +7.76%  '
$ws.Cells.Item(44, 5).Value = '  +9.45%  '
$ws.Cells.Item(45, 4).Value = '''9.29'
$ws.Cells.Item(45, 5).Value = '  +9.24%  '
$ws.Cells.Item(46, 4).Value = '''62.13'
$ws.Cells.Item(46, 5).Value = '  +1.06%  '
$ws.Cells.Item(47, 4).Value = '''4.90'
$ws.Cells.Item(47, 5).Value = '  -4.09%  '
$ws.Cells.Item(48, 5).Value = '  +3.75%  '
$ws.Cells.Item(49, 2).Value = 'ARBITRUM'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(49, 4).Value = '''1.18'
$ws.Cells.Item(49, 5).Value = '  +3.17%  '
$ws.Cells.Item(50, 2).Value = 'BinanceUSD'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(50, 4).Value = '''1.00'
$ws.Cells.Item(50, 5).Value = '  +0.05%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '''97.19'
$ws.Cells.Item(51, 5).Value = '  +6.78%  '
